$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DIEGO's balance (row 5, column C "Saldo") from 32314.16 to 33000
$ws.Cells.Item(5, 3).Value = 33000

# Remove the two rows that were deleted in the source edit.
# Delete from the bottom up so earlier row indices stay valid.
# Row 9 = RAFAELA / 002064834 (balance 1988.37)
$ws.Rows.Item(9).Delete()
# Row 6 = E3 / 004267976 (balance 28640.83)
$ws.Rows.Item(6).Delete()
